# Rework the counters_summary sheet: new set of dimension columns, each
# with a paired "... SCORE" column, replacing the old ad-hoc column layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (column A "attribute" stays as-is)
$headers = @(
    "COMPLETENESSMANDATORY",
    "COMPLETENESSMANDATORY SCORE",
    "COMPLETENESSOPTIONAL",
    "COMPLETENESSOPTIONAL SCORE",
    "PRECISION",
    "PRECISION SCORE",
    "BUSINESSRULECOMPLIANCE",
    "BUSINESSRULECOMPLIANCE SCORE",
    "METADATACOMPLIANCE",
    "METADATACOMPLIANCE SCORE",
    "UNIQUENESS",
    "UNIQUENESS SCORE",
    "NONREDUNDANCY",
    "NONREDUNDANCY SCORE",
    "SEMANTICCONSISTENCY",
    "SEMANTICCONSISTENCY SCORE",
    "VALUECONSISTENCY",
    "VALUECONSISTENCY SCORE",
    "FORMATCONSISTENCY",
    "FORMATCONSISTENCY SCORE"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Find last used data row (attribute names live in column A)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Reset every score/measure cell for the new column range (B..U) to 0 for
# each existing data row.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le 21; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}
